# Tech limit / tech freeze future scenario
# Re-allocates part of the residual delta-cost trend (deltaC_Res_Ops, col K)
# into the structural delta-cost trend (deltaC_Structural_Ops, col G) for the
# operational scenario, row by row (years 1960-2020), while leaving the
# operational total (deltaC_Tot_Ops, col L = G+H+I+J+K) unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: row number -> @(new deltaC_Structural_Ops (G), new deltaC_Res_Ops (K))
$techScenario = @{
    2 = @(2.664992669550213, -3.044232424142641)
    3 = @(5.22944814093343, 0.008781852903227971)
    4 = @(7.701324314756554, 0.529162567745006)
    5 = @(10.08502655826333, -1.393737462944404)
    6 = @(12.3824108185826, -6.582795513075803)
    7 = @(14.59355353798953, -7.202028896358314)
    8 = @(16.71734054847584, -10.60487122550825)
    9 = @(18.75191837698277, -9.001206383730835)
    10 = @(20.69503881460975, -1.050647309108927)
    11 = @(22.54432173840143, 3.279922686611684)
    12 = @(24.29745494400631, -1.635340019385231)
    13 = @(25.95234446962981, 1.548197332326826)
    14 = @(27.50722739444926, -3.563079834944676)
    15 = @(28.96075377193128, -3.036326780622455)
    16 = @(30.3120461043778, -4.850337260100135)
    17 = @(31.5607398441365, -2.781746060624252)
    18 = @(32.70700815622536, -3.135847663060898)
    19 = @(33.75157581675575, -2.676769608168425)
    20 = @(34.69572345247781, -8.954580260379982)
    21 = @(35.54128298783034, -10.25115260178631)
    22 = @(36.2906284256927, -2.476404760556565)
    23 = @(36.94665876890706, -1.894768268810874)
    24 = @(37.51278097172949, -0.03841230759252845)
    25 = @(37.99288479471285, 0.2733501000712391)
    26 = @(38.39131892043434, 0.5674213424777861)
    27 = @(38.71286252940088, 0.02625986414867043)
    28 = @(38.96269685331959, 3.083091979982411)
    29 = @(39.14637564117176, -0.5082175687428219)
    30 = @(39.26979423973023, -0.2853962880410066)
    31 = @(39.33915982804758, -0.1925172440725333)
    32 = @(39.36096207950109, 1.653299560109701)
    33 = @(39.34194435857267, 5.416849539083154)
    34 = @(39.28907739495488, 7.584990792866384)
    35 = @(39.20953391717236, 11.08564651495203)
    36 = @(39.11066735850467, 8.010344797365965)
    37 = @(38.99999225307309, 6.998918502473568)
    38 = @(38.88516981348979, 4.481702914425053)
    39 = @(38.77399697719627, 3.124660224622911)
    40 = @(38.67440060583528, 4.857910953413057)
    41 = @(38.59443687368142, 3.920862048942016)
    42 = @(38.54229680382326, -0.3348749073475474)
    43 = @(38.52631843076895, 4.473823206108992)
    44 = @(38.55500553267243, 0.01674689534563356)
    45 = @(38.63705233192732, 0.03868216963427429)
    46 = @(38.78137783272454, -3.564465671286293)
    47 = @(38.99716536776136, -6.417030641841002)
    48 = @(39.2939107908675, -7.555404318087341)
    49 = @(39.68147723019083, -8.809888048050126)
    50 = @(40.1701566951564, -5.700446577307943)
    51 = @(40.77073842495055, -5.790901603722958)
    52 = @(41.49458234716073, -7.603550163739243)
    53 = @(42.35369782390887, -5.177955826106512)
    54 = @(43.36082611602388, -4.897993431460492)
    55 = @(44.52952543831268, -3.699918079014267)
    56 = @(45.87425774419418, -1.15360274710666)
    57 = @(47.41047671287962, 1.119302996120076)
    58 = @(49.15471372923957, 5.105825089846185)
    59 = @(51.12466315781803, 6.615089869845292)
    60 = @(53.33926304455639, 11.52251000415518)
    61 = @(55.8187734274934, 15.30326424878677)
    62 = @(58.58484741858675, 76.46869111733241)
}

foreach ($row in $techScenario.Keys) {
    $values = $techScenario[$row]
    $ws.Cells.Item($row, 7).Value = $values[0]   # G: deltaC_Structural_Ops
    $ws.Cells.Item($row, 11).Value = $values[1]  # K: deltaC_Res_Ops
}
